$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = "{OKRs}"
$ws.Range("F11").Value = "{SemesterGrade}"

$ws.Range("H15").Select() | Out-Null
